$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Records" table (column J, header at J3="Records") currently lists:
#   J4=Id, J5=Date, J6=Obs, J7=ExamFK, J8=AttendFK, J9=EnrollFK
# Remove the "Obs" field, shifting the following FK fields up one row.
$ws.Range("J6").Value = "ExamFK"
$ws.Range("J7").Value = "AttendFK"
$ws.Range("J8").Value = "EnrollFK"
$ws.Range("J9").ClearContents()

# "Students" table (column L) gains a new field "ProfileId" in the next
# free row.
$ws.Range("L9").Value = "ProfileId"

# Update the active selection to reflect where the editor was last
# working, matching the saved workbook state.
$ws.Range("K24").Select()
